$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: a second "jiji" patient record, mirroring the layout/format of
# row 2 but with its own date values and a couple of columns shifted
# (no H4 "耐药" value; it now lands in I4 instead), plus new lookup strings.

# A4: plain text (same shared string as A2 "jiji")
$ws.Range("A4").Value = "jiji"

# B4: numeric "33" styled with the text-like "@" number format (same as B2)
$ws.Range("B4").Value = 33
$ws.Range("B4").NumberFormat = "@"

# C4: date (serial 42758 == 2017-01-23), same date style as C2
$ws.Range("C2").Copy($ws.Range("C4"))
$ws.Range("C4").Value = 42758

# D4: numeric "23" with the same 0.00 style as D2
$ws.Range("D2").Copy($ws.Range("D4"))
$ws.Range("D4").Value = 23

# E4: date (serial 42361 == 2015-12-23), same date style as E2
$ws.Range("E2").Copy($ws.Range("E4"))
$ws.Range("E4").Value = 42361

# F4: same "3,4" text used in F2
$ws.Range("F4").Value = "3,4"

# G4: new text "牛逼2"
$ws.Range("G4").Value = "牛逼2"

# I4: "耐药" (same string as H2), now in column I for this row
$ws.Range("I4").Value = "耐药"

# J4: new text "鸡巴"
$ws.Range("J4").Value = "鸡巴"

# Update the active selection to match the new edit location
$ws.Range("F4").Select()
